# Update the LR-pairs table for Cntn2-Cntn1 with the new TPM-derived values.
# The data was recomputed: the "Neutrophils" sending cluster was relabeled to
# "Inflammatory-Mac", and all rows whose Target cluster was "MuSCs" were
# dropped (only "FAPs" target rows remain), leaving 3 data rows instead of 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that targeted "MuSCs" (old rows 3, 5 and 7).
# Deleted bottom-up so row numbers of rows not yet processed stay valid.
$ws.Rows("7:7").Delete()
$ws.Rows("5:5").Delete()
$ws.Rows("3:3").Delete()

# After the deletions the sheet has 3 data rows left (old rows 2, 4, 6 -
# now rows 2, 3, 4), in order ECs / MuSCs / Neutrophils. Overwrite every
# cell with the final values from the recomputed TPM data so that both
# the remaining row order (ECs, Inflammatory-Mac, MuSCs) and all of the
# numbers are correct, regardless of how the deletes re-numbered rows.

# Row 2: ECs -> Cntn2 -> Cntn1 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cntn2"
$ws.Range("C2").Value = "Cntn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.0700585
$ws.Range("H2").Value = 0.140117
$ws.Range("I2").Value = 0.4535220567529688
$ws.Range("J2").Value = 0.4377808049690373
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1145166666666667
$ws.Range("N2").Value = 0.34355
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.008022865891666666
$ws.Range("R2").Value = 0.04813719535
$ws.Range("S2").Value = 0.4535220567529688
$ws.Range("T2").Value = 0.4377808049690373

# Row 3: Inflammatory-Mac (formerly "Neutrophils") -> Cntn2 -> Cntn1 -> FAPs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Cntn2"
$ws.Range("C3").Value = "Cntn1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.011109
$ws.Range("H3").Value = 0.033327
$ws.Range("I3").Value = 0.07191385097409637
$ws.Range("J3").Value = 0.1041267004517875
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1145166666666667
$ws.Range("N3").Value = 0.34355
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.00127216565
$ws.Range("R3").Value = 0.01144949085
$ws.Range("S3").Value = 0.07191385097409637
$ws.Range("T3").Value = 0.1041267004517875

# Row 4: MuSCs -> Cntn2 -> Cntn1 -> FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cntn2"
$ws.Range("C4").Value = "Cntn1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.073309
$ws.Range("H4").Value = 0.146618
$ws.Range("I4").Value = 0.4745640922729347
$ws.Range("J4").Value = 0.4580924945791753
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1145166666666667
$ws.Range("N4").Value = 0.34355
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.008395102316666667
$ws.Range("R4").Value = 0.0503706139
$ws.Range("S4").Value = 0.4745640922729347
$ws.Range("T4").Value = 0.4580924945791753
